# Auto-generated edit script for cryptos.xlsx update
# Applies the "Updated cryptos list" commit changes to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cell updates (values that are not ambiguous with numbers) ---
$ws.Range('D2').Value = '64.962.73'
$ws.Range('E2').Value = '  +3.66%  '
$ws.Range('D3').Value = '2.536.84'
$ws.Range('E3').Value = '  +3.04%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('E6').Value = '  +3.90%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +1.27%  '
$ws.Range('D9').Value = '2.538.20'
$ws.Range('E9').Value = '  +3.05%  '
$ws.Range('E10').Value = '  +2.11%  '
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('E12').Value = '  +0.57%  '
$ws.Range('E13').Value = '  +0.68%  '
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('E15').Value = '  +2.76%  '
$ws.Range('D16').Value = '2.997.22'
$ws.Range('E16').Value = '  +3.03%  '
$ws.Range('D17').Value = '64.723.81'
$ws.Range('E17').Value = '  +3.37%  '
$ws.Range('D18').Value = '2.546.65'
$ws.Range('E18').Value = '  +3.40%  '
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('E21').Value = '  +3.62%  '
$ws.Range('E22').Value = '  +1.38%  '
$ws.Range('E23').Value = '  +2.98%  '
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('E28').Value = '  +7.90%  '
$ws.Range('D29').Value = '2.673.62'
$ws.Range('E29').Value = '  +3.58%  '
$ws.Range('E30').Value = '  +4.71%  '
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('E32').Value = '  +2.00%  '
$ws.Range('E33').Value = '  +2.67%  '
$ws.Range('E34').Value = '  +2.49%  '
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('E36').Value = '  +3.40%  '
$ws.Range('E37').Value = '  +2.74%  '
$ws.Range('E38').Value = '  +6.20%  '
$ws.Range('E39').Value = '  +2.50%  '
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('E42').Value = '  +1.84%  '
$ws.Range('E43').Value = '  +5.77%  '
$ws.Range('E44').Value = '  +5.69%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').Value = '0.0₆0300'
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('E47').Value = '  +2.31%  '
$ws.Range('E48').Value = '  +2.73%  '
$ws.Range('E50').Value = '  +2.79%  '
$ws.Range('E51').Value = '  +2.17%  '

# --- Cell updates that look numeric and must be forced to remain text ---
# (format each cell as text, assign the value, then clear the temporary
#  number format so the cell style index matches the original unstyled cells)
$forcedTextCells = @(
    'D5',
    'D6',
    'D12',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D25',
    'D26',
    'D27',
    'D31',
    'D32',
    'D36',
    'D38',
    'D39',
    'D43',
    'D44',
    'D47',
    'D49',
    'D51'
)
foreach ($cellRef in $forcedTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}
$ws.Range('D5').Value = '581.10'
$ws.Range('D6').Value = '152.81'
$ws.Range('D12').Value = '5.30'
$ws.Range('D19').Value = '8.01'
$ws.Range('D20').Value = '11.01'
$ws.Range('D21').Value = '4.28'
$ws.Range('D22').Value = '329.71'
$ws.Range('D23').Value = '2.24'
$ws.Range('D25').Value = '10.24'
$ws.Range('D26').Value = '65.89'
$ws.Range('D27').Value = '643.61'
$ws.Range('D31').Value = '0.993'
$ws.Range('D32').Value = '8.08'
$ws.Range('D36').Value = '1.58'
$ws.Range('D38').Value = '5.64'
$ws.Range('D39').Value = '154.93'
$ws.Range('D43').Value = '1.83'
$ws.Range('D44').Value = '161.85'
$ws.Range('D47').Value = '15.65'
$ws.Range('D49').Value = '21.54'
$ws.Range('D51').Value = '0.0518'
foreach ($cellRef in $forcedTextCells) {
    $ws.Range($cellRef).ClearFormats()
}

